$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "listsProd"
$ws.Range("C1").Value = "listsFC"
$ws.Range("B2").Value = "Lists/voicingProd.xlsx"
$ws.Range("C2").Value = "Lists/voicingFC.xlsx"
$ws.Range("B3").Value = "Lists/devoicingProd.xlsx"
$ws.Range("C3").Value = "Lists/devoicingFC.xlsx"

$ws.Range("E6").Select()
